$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so values can be written, then
# re-apply the same protection settings afterward.
$ws.Unprotect("D382")

# Update the confidential footer note date (A9): 2021-04-27 -> 2021-04-28
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for rows 2-6
$ws.Range("D2").Value = 0.2517662468256127
$ws.Range("E2").Value = 0.0027650877614811

$ws.Range("D3").Value = 0.250078411121557
$ws.Range("E3").Value = 0.003910614525139522

$ws.Range("D4").Value = 0.2477546960721267
$ws.Range("E4").Value = -0.002056807051909804

$ws.Range("D5").Value = 0.2504006459807035
$ws.Range("E5").Value = 0.00917431192660545

$ws.Range("E6").Value = 0.003461786061632832

# Restore sheet protection to match the original workbook state
$ws.Protect("D382")
